$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts existing "Tipo" column from D to E)
$ws.Range("D1").EntireColumn.Insert()

# Copy formatting from the neighboring header cell (C1) onto the new header cell (D1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set header text for new MAE column
$ws.Range("D1").Value = "MAE"

# Set the MAE value for row 2
$ws.Range("D2").Value = 0.3447233267553575

$ws.Dimension = "A1:E2"
